# Team07Report.xlsx edit script
# Marks the "Marriage before death" (US05) and "Divorce before death" (US06)
# stories/tasks as Done, records actual size/time for the completed work,
# reduces the Burndown remaining-points count for the latest sprint data
# point, and leaves the various sheets scrolled/selected where the author
# last left them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Backlog: mark the corresponding backlog rows for US05 / US06 as Done.
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets("Backlog")
$backlog.Activate()
$backlog.Range("E13").Value = "Done"
$backlog.Range("E16").Value = "Done"
$backlog.Range("E17").Value = "Done"
$backlog.Range("E18").Select()

# ---------------------------------------------------------------------
# Sprint1: no data changes, only the remembered scroll position moved.
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets("Sprint1")
$sprint1.Activate()
$excel.ActiveWindow.ScrollRow = 8
$sprint1.Range("H22").Select()

# ---------------------------------------------------------------------
# Sprint2: US05 (row 7) and US06 (row 12), plus their sub-tasks
# (rows 8-11 and 13-16), move from "Coding" to "Done". US05 and US06
# also get their actual size/time filled in, and US05 is flagged
# completed.
# ---------------------------------------------------------------------
$sprint2 = $wb.Worksheets("Sprint2")
$sprint2.Activate()

$sprint2.Range("D7").Value = "Done"
$sprint2.Range("G7").Value = 20
$sprint2.Range("H7").Value = 20
$sprint2.Range("I7").Value = "yes"

$sprint2.Range("D8").Value = "Done"
$sprint2.Range("D9").Value = "Done"
$sprint2.Range("D10").Value = "Done"
$sprint2.Range("D11").Value = "Done"

$sprint2.Range("D12").Value = "Done"
$sprint2.Range("G12").Value = 20
$sprint2.Range("H12").Value = 20

$sprint2.Range("D13").Value = "Done"
$sprint2.Range("D14").Value = "Done"
$sprint2.Range("D15").Value = "Done"
$sprint2.Range("D16").Value = "Done"

$sprint2.Range("I7").Select()

# ---------------------------------------------------------------------
# Stories: record the owner ("es") for the US05 / US06 rows.
# ---------------------------------------------------------------------
$stories = $wb.Worksheets("Stories")
$stories.Activate()
$stories.Range("D6").Value = "es"
$stories.Range("D7").Value = "es"
$excel.ActiveWindow.ScrollRow = 3
$stories.Range("D7").Select()

# ---------------------------------------------------------------------
# Burndown: the sprint's remaining points dropped from 27 to 25, which
# ripples through the delta and the (now-finished) Sprint2 totals.
# This is also the sheet left active/selected at the end.
# ---------------------------------------------------------------------
$burndown = $wb.Worksheets("Burndown")
$burndown.Activate()
$burndown.Range("C4").Value = 25
$burndown.Range("D8").Select()
